# Covid-19 KSA analysis v2 final draft.pptx
#
# 1) Refresh the cached "datetimeFigureOut" date field (3/6/2024 -> 3/27/2024)
#    on the Slide Master and on every Slide Layout's Date placeholder.
# 2) Fix two typos on slide 13 ("locdown" -> "lockdown", "untiol" -> "until"),
#    which also collapses the previously split runs into a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder refresh
# ---------------------------------------------------------------------------
$newDate = "3/27/2024"

function Update-DateField($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Type -eq 14 -and $sh.HasTextFrame) {
            $ph = $sh.PlaceholderFormat
            if ($ph.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateField $layouts.Item($i).Shapes
}

# ---------------------------------------------------------------------------
# 2) Typo fixes on slide 13
#    A direct Text= on a range spanning several runs only patches the
#    characters that actually differ, leaving the old run split (and its
#    err="1" spell-check marks) in place. Blow the old runs away first with
#    a throwaway placeholder string (so nothing lines up with the final
#    text) and then write the corrected sentence - that collapses the
#    paragraph back down to a single run using the first run's formatting.
# ---------------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)

$rect4 = $slide13.Shapes.Item(4)
$rect4.TextFrame.TextRange.Text = "#"
$rect4.TextFrame.TextRange.Text = "No quarantine relaxation on the cards as yet lockdown to continue until 14 May"

$rect11 = $slide13.Shapes.Item(8)
$rect11.TextFrame.TextRange.Text = "#"
$rect11.TextFrame.TextRange.Text = "No quarantine relaxation on the cards as yet and lockdown to remain until 11 May"
